$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 685, shifting existing rows 685-727 down to 686-728
$ws.Rows.Item(685).Insert()

# Populate the newly inserted row 685 with the new weekly price record
$ws.Cells.Item(685, 1).Value  = 6
$ws.Cells.Item(685, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(685, 3).Value  = "Metropolitana"
$ws.Cells.Item(685, 4).Value  = 45041
$ws.Cells.Item(685, 5).Value  = 13
$ws.Cells.Item(685, 6).Value  = 100112044
$ws.Cells.Item(685, 7).Value  = "Perejil"
$ws.Cells.Item(685, 8).Value  = "Sin especificar"
$ws.Cells.Item(685, 9).Value  = "Primera"
$ws.Cells.Item(685, 10).Value = 290
$ws.Cells.Item(685, 11).Value = 12000
$ws.Cells.Item(685, 12).Value = 13000
$ws.Cells.Item(685, 13).Value = 12448
$ws.Cells.Item(685, 14).Value = "$/docena de atados"
$ws.Cells.Item(685, 15).Value = "Región Metropolitana"
$ws.Cells.Item(685, 16).Value = 4149
$ws.Cells.Item(685, 17).Value = 3
$ws.Cells.Item(685, 18).Value = "Hortaliza"
